$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Paris 2023 Contenders Sticker Capsule'
$ws.Range("B1").Value = '$0.28 USD'
